$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the changed cells to text format so numeric-looking values
# (e.g. "1.002") are stored as text, matching the source data, then
# clear the formatting afterwards so no extra style index is left
# behind on the cells (the original cells carry no explicit style).
$changedRange = $ws.Range("D2:E51")
$changedRange.NumberFormat = "@"

$ws.Range("D2").Value = "22.560.65"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.578.15"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").Value = "288.75"
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("E7").Value = "  -1.05%  "
$ws.Range("D8").Value = "48.52"
$ws.Range("E8").Value = "  -2.91%  "
$ws.Range("D9").Value = "0.3357"
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("D10").Value = "1.145"
$ws.Range("E10").Value = "  +0.27%  "
$ws.Range("D11").Value = "0.07487"
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("D14").Value = "6.013"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").Value = "6.963"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").Value = "1.584.90"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").Value = "0.00001118"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "88.75"
$ws.Range("E18").Value = "  -2.13%  "
$ws.Range("D19").Value = "0.06760"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "6.428"
$ws.Range("E20").Value = "  +2.29%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("E22").Value = "  +1.23%  "
$ws.Range("D23").Value = "12.19"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "22.563.09"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("D25").Value = "2.409"
$ws.Range("E25").Value = "  +1.94%  "
$ws.Range("D26").Value = "2.610"
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("D27").Value = "152.97"
$ws.Range("E27").Value = "  +1.89%  "
$ws.Range("D28").Value = "19.66"
$ws.Range("E28").Value = "  -1.89%  "
$ws.Range("D29").Value = "5.014"
$ws.Range("E29").Value = "  -0.86%  "
$ws.Range("D30").Value = "124.53"
$ws.Range("E30").Value = "  -0.40%  "
$ws.Range("D31").Value = "1.758.45"
$ws.Range("E31").Value = "  +0.69%  "
$ws.Range("D32").Value = "1.079"
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("D33").Value = "6.210"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").Value = "2.005"
$ws.Range("E34").Value = "  -0.23%  "
$ws.Range("D35").Value = "9.805"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "0.08329"
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("D37").Value = "0.02462"
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("D38").Value = "0.2273"
$ws.Range("D39").Value = "5.493"
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("D40").Value = "0.06461"
$ws.Range("E40").Value = "  -1.37%  "
$ws.Range("D41").Value = "1.303"
$ws.Range("E41").Value = "  -2.84%  "
$ws.Range("D42").Value = "11.46"
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("D43").Value = "0.6352"
$ws.Range("E43").Value = "  +1.94%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "14.01"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.6184"
$ws.Range("E45").Value = "  +5.51%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "3.773"
$ws.Range("E46").Value = "  -1.05%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "2.069"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "125.96"
$ws.Range("E48").Value = "  -2.54%  "
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").Value = "1.224"
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.07285"
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "77.31"
$ws.Range("E51").Value = "  +0.92%  "

$changedRange.ClearFormats()
